# Generate Report for Handoff
#
# Two new source files go through handoff:
#   2ddde282-c7ce-4e0b-9f58-529427073233.md
#   6286292a-6a47-4e45-b807-843655502b2f.md
#
# They are inserted (in that order) ahead of the ".localization-config"
# bookkeeping row, which is pushed down, on the "Overview" sheet and on
# each per-locale detail sheet ("zh-cn", "de-de").
#
# NOTE: this runtime's Hyperlinks.Add only replaces a hyperlink that was
# itself added earlier in the same session - a hyperlink that came from
# the loaded file sticks around even after the cell's contents/hyperlink
# are reassigned, producing duplicate <hyperlink> entries for the same
# ref. So every sheet's hyperlinks are cleared up front and the full set
# (unchanged rows included) is re-added in ref order, which also keeps
# the r:id numbering stable/sequential.

$wb = $excel.ActiveWorkbook

$mdBase      = "https://github.com/OpenLocalizationTest/oltest/blob/17e6234b63dd3fe3b91908e55183a990462be170/e2e/"
$configUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/17e6234b63dd3fe3b91908e55183a990462be170/.localization-config"
$zhcnXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2fb5a5e1c86369d9b7aedfeff1fbfde50b52110c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/"
$dedeXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e6cdb4e9fa4723862aaac2804a9a14cc09dfefc9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/"

$mdFile1 = "afe93a64-a777-4d97-9a4f-72fa108740f4.md"
$mdFile2 = "b96deeef-3b69-46dc-9968-44c0cb3c132c.md"
$mdFile3 = "2ddde282-c7ce-4e0b-9f58-529427073233.md"
$mdFile4 = "6286292a-6a47-4e45-b807-843655502b2f.md"
$configFile = ".localization-config"

$file3Hash = "a4cc2a326851d19300fcbb1754cf8fc4f63f16a7"
$file3ZhCnXlf = "2ddde282-c7ce-4e0b-9f58-529427073233.$file3Hash.zh-cn.xlf"
$file3DeDeXlf = "2ddde282-c7ce-4e0b-9f58-529427073233.$file3Hash.de-de.xlf"

$file4Hash = "3296402e1cb7304c1b31eeed076b428c6057e7ae"
$file4ZhCnXlf = "6286292a-6a47-4e45-b807-843655502b2f.$file4Hash.zh-cn.xlf"
$file4DeDeXlf = "6286292a-6a47-4e45-b807-843655502b2f.$file4Hash.de-de.xlf"

$file1ZhCnXlf = "afe93a64-a777-4d97-9a4f-72fa108740f4.c96f027e0f897c5a3f6fd602f610f0137e4e2108.zh-cn.xlf"
$file1DeDeXlf = "afe93a64-a777-4d97-9a4f-72fa108740f4.c96f027e0f897c5a3f6fd602f610f0137e4e2108.de-de.xlf"
$file2ZhCnXlf = "b96deeef-3b69-46dc-9968-44c0cb3c132c.896969fcfb36782802e7faa50f9957d30c65a5b3.zh-cn.xlf"
$file2DeDeXlf = "b96deeef-3b69-46dc-9968-44c0cb3c132c.896969fcfb36782802e7faa50f9957d30c65a5b3.de-de.xlf"

$inTranslation = "In Translation"
$readyStatus = "Ready for handoff"
$notLocalizedStatus = "Not to be localized"
$epoch = "0001-01-01 00:00:00"
$ignoredReason = "Ignored"
$includeReason = "Include"

$handoffDt12 = "2016-03-02 14:11:56"
$handoffDt13 = "2016-03-02 14:12:18"
$handoffDt34zhcn = "2016-03-02 14:13:41"
$handoffDt34dede = "2016-03-02 14:13:53"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $mdFile3
$wsOverview.Range("B4").Value = $readyStatus
$wsOverview.Range("C4").Value = $readyStatus

$wsOverview.Range("A5").Value = $mdFile4
$wsOverview.Range("B5").Value = $readyStatus
$wsOverview.Range("C5").Value = $readyStatus

$wsOverview.Range("A6").Value = $configFile
$wsOverview.Range("B6").Value = $notLocalizedStatus
$wsOverview.Range("C6").Value = $notLocalizedStatus

$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdBase + $mdFile1, "", "", $mdFile1) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdBase + $mdFile2, "", "", $mdFile2) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $mdBase + $mdFile3, "", "", $mdFile3) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), $mdBase + $mdFile4, "", "", $mdFile4) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), $configUrl, "", "", $configFile) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn": detailed handoff/handback tracking for the zh-cn locale
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("B2").Value = $inTranslation
$wsZhCn.Range("D2").Value = $handoffDt12
$wsZhCn.Range("G2").Value = $epoch
$wsZhCn.Range("H2").Value = $includeReason

$wsZhCn.Range("B3").Value = $inTranslation
$wsZhCn.Range("D3").Value = $handoffDt12
$wsZhCn.Range("G3").Value = $epoch
$wsZhCn.Range("H3").Value = $includeReason

$wsZhCn.Range("A4").Value = $mdFile3
$wsZhCn.Range("B4").Value = $readyStatus
$wsZhCn.Range("D4").Value = $handoffDt34zhcn
$wsZhCn.Range("G4").Value = $epoch
$wsZhCn.Range("H4").Value = $includeReason

$wsZhCn.Range("A5").Value = $mdFile4
$wsZhCn.Range("B5").Value = $readyStatus
$wsZhCn.Range("D5").Value = $handoffDt34zhcn
$wsZhCn.Range("G5").Value = $epoch
$wsZhCn.Range("H5").Value = $includeReason

$wsZhCn.Range("A6").Value = $configFile
$wsZhCn.Range("B6").Value = $notLocalizedStatus
$wsZhCn.Range("D6").Value = $epoch
$wsZhCn.Range("G6").Value = $epoch
$wsZhCn.Range("H6").Value = $ignoredReason

$wsZhCn.Range("A1").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdBase + $mdFile1, "", "", $mdFile1) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), $zhcnXlfBase + $file1ZhCnXlf, "", "", $file1ZhCnXlf) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $mdBase + $mdFile2, "", "", $mdFile2) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C3"), $zhcnXlfBase + $file2ZhCnXlf, "", "", $file2ZhCnXlf) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), $mdBase + $mdFile3, "", "", $mdFile3) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C4"), $zhcnXlfBase + $file3ZhCnXlf, "", "", $file3ZhCnXlf) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), $mdBase + $mdFile4, "", "", $mdFile4) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C5"), $zhcnXlfBase + $file4ZhCnXlf, "", "", $file4ZhCnXlf) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), $configUrl, "", "", $configFile) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de": detailed handoff/handback tracking for the de-de locale
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("B2").Value = $inTranslation
$wsDeDe.Range("D2").Value = $handoffDt13
$wsDeDe.Range("G2").Value = $epoch
$wsDeDe.Range("H2").Value = $includeReason

$wsDeDe.Range("B3").Value = $inTranslation
$wsDeDe.Range("D3").Value = $handoffDt13
$wsDeDe.Range("G3").Value = $epoch
$wsDeDe.Range("H3").Value = $includeReason

$wsDeDe.Range("A4").Value = $mdFile3
$wsDeDe.Range("B4").Value = $readyStatus
$wsDeDe.Range("D4").Value = $handoffDt34dede
$wsDeDe.Range("G4").Value = $epoch
$wsDeDe.Range("H4").Value = $includeReason

$wsDeDe.Range("A5").Value = $mdFile4
$wsDeDe.Range("B5").Value = $readyStatus
$wsDeDe.Range("D5").Value = $handoffDt34dede
$wsDeDe.Range("G5").Value = $epoch
$wsDeDe.Range("H5").Value = $includeReason

$wsDeDe.Range("A6").Value = $configFile
$wsDeDe.Range("B6").Value = $notLocalizedStatus
$wsDeDe.Range("D6").Value = $epoch
$wsDeDe.Range("G6").Value = $epoch
$wsDeDe.Range("H6").Value = $ignoredReason

$wsDeDe.Range("A1").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdBase + $mdFile1, "", "", $mdFile1) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), $dedeXlfBase + $file1DeDeXlf, "", "", $file1DeDeXlf) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $mdBase + $mdFile2, "", "", $mdFile2) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C3"), $dedeXlfBase + $file2DeDeXlf, "", "", $file2DeDeXlf) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), $mdBase + $mdFile3, "", "", $mdFile3) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C4"), $dedeXlfBase + $file3DeDeXlf, "", "", $file3DeDeXlf) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), $mdBase + $mdFile4, "", "", $mdFile4) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C5"), $dedeXlfBase + $file4DeDeXlf, "", "", $file4DeDeXlf) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), $configUrl, "", "", $configFile) | Out-Null

Write-Output "Generated handoff report rows for $mdFile3 and $mdFile4"
